$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-16 Monday" "2026-02-17 Tuesday"

Replace-Text "238÷7=" "257÷6="
Replace-Text "248÷9=" "571÷2="
Replace-Text "177÷4=" "751÷3="
Replace-Text "111÷6=" "858÷4="
Replace-Text "776÷5=" "207÷8="

Replace-Text "767÷4=" "667÷2="
Replace-Text "894÷2=" "791÷8="
Replace-Text "567÷3=" "535÷4="
Replace-Text "881÷2=" "604÷2="
Replace-Text "460÷3=" "688÷2="

Replace-Text "148÷2=" "289÷2="
Replace-Text "666÷3=" "585÷8="
Replace-Text "725÷3=" "278÷5="
Replace-Text "302÷2=" "928÷5="
Replace-Text "978÷8=" "811÷3="

Replace-Text "609÷8=" "268÷5="
Replace-Text "915÷8=" "297÷8="
Replace-Text "171÷9=" "486÷8="
Replace-Text "592÷7=" "595÷8="
Replace-Text "492÷5=" "382÷5="

Replace-Text "102÷8=" "693÷8="
Replace-Text "719÷2=" "970÷8="
Replace-Text "411÷3=" "428÷6="
Replace-Text "720÷5=" "850÷6="
Replace-Text "623÷4=" "248÷6="
